$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19, pushing existing rows 19-86 down to 20-87.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new weekly record.
$ws.Range("A19").Value = 11
$ws.Range("B19").Value = "Vega Monumental Concepción"
$ws.Range("C19").Value = "Bíobío"
$ws.Range("D19").Value = 44565
$ws.Range("E19").Value = 8
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100108
$ws.Range("H19").Value = "Tropicales y subtropicales"
$ws.Range("I19").Value = 100108002
$ws.Range("J19").Value = "Mango"
$ws.Range("K19").Value = "Sin especificar"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 140
$ws.Range("N19").Value = 6500
$ws.Range("O19").Value = 7000
$ws.Range("P19").Value = 6786
$ws.Range("Q19").Value = '$/bandeja 4 kilos'
$ws.Range("R19").Value = "Perú"
$ws.Range("S19").Value = 1696
$ws.Range("T19").Value = 4

Write-Host "done"
